# Timesheet update: "Updated log until 9th march"
#
# - Row 21 ("KTN Implementation user and server") gets its label extended
#   to "...login" and a Torsdag (K) hour entry of 4.
# - Row 5 ("Documentation") gains a Torsdag (K) hour entry of 1.
# - Two brand-new activity rows are inserted after row 22
#   ("MMI Report" and "Java fx installation and understanding…"),
#   pushing everything below down by two rows.
# - The now-shifted blank row right after those two gets filled in with
#   the "DB: Implementation of database in java." entry.
# - The weekly SUM() formulas in row 3 automatically pick up the two
#   extra rows because the inserts happen inside their summed range.
# - Selection cursor moves to M10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the (currently blank) row 23. This shifts
# row 23 onward down by two and extends every SUM(...) range that spans
# across the insertion point (rows 6:30 -> 6:32, 5:30 -> 5:32, etc.)
# while preserving the shared-formula grouping in row 3.
$ws.Rows.Item(23).Resize(2).Insert()

# Row 5: "Documentation" - add Thursday hours.
$ws.Range("K5").Value = 1

# Row 21: relabel and add Thursday hours.
$ws.Range("A21").Value = "KTN Implementation user and server login"
$ws.Range("K21").Value = 4

# Newly inserted row 23: "MMI Report".
$ws.Range("A23").Value = "MMI Report"
$ws.Range("K23").Value = 14

# Newly inserted row 24: "Java fx installation and understanding...".
$ws.Range("A24").Value = "Java fx installation and understanding…"
$ws.Range("K24").Value = 11

# Row 25 (formerly the blank row 23, now shifted down by the insert):
# "DB: Implementation of database in java."
$ws.Range("A25").Value = "DB: Implementation of database in java. "
$ws.Range("K25").Value = 7

# Move the active selection, matching the saved view state.
$ws.Range("M10").Select() | Out-Null
